$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.206.37'
$ws.Range("E2").Value = '  +1.00%  '
$ws.Range("D3").Value = '1.802.05'
$ws.Range("E3").Value = '  +2.60%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.24%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '335.94'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '1.000'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +0.09%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4609'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +20.68%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3712'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +10.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '45.03'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +0.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.151'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07611'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +5.70%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.002'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.10%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.37'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.342'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.58%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.425'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.62%  '
$ws.Range("D16").Value = '1.799.07'
$ws.Range("E16").Value = '  +2.41%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001103'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.57%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.06720'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.92%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '82.15'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.10%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.000'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +0.01%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.47'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.68%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.404'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.02%  '
$ws.Range("D23").Value = '28.177.87'
$ws.Range("E23").Value = '  +0.79%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '11.90'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.59%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.412'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.88%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '20.71'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.42%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.381'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +3.21%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.44'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").Value = '2.004.09'
$ws.Range("E29").Value = '  +2.30%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '133.89'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.60%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.256'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.03%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.037'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.47%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.09576'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +8.98%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.876'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.99%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.2223'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.45%  '
$ws.Range("B36").Value = 'Hedera'
$ws.Range("C36").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06367'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.62%  '
$ws.Range("B37").Value = 'VeChain'
$ws.Range("C37").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02357'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +2.33%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '12.06'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.05%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.253'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.70%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.6657'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +0.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.520'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +3.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.233'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.088'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.05%  '
$ws.Range("B44").Value = 'Frax'
$ws.Range("C44").Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.9998'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +0.05%  '
$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.10'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.6101'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.39%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.828'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.21%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '130.72'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +3.92%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.052'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +2.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.180'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.84%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.07159'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.34%  '
